$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Times")

$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = "Arran Katoko FC"
